# Update "想去人数" (want-to-go count) figures across sheets to match the
# newly generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 7
$ws1.Range("F6").Value = 159
$ws1.Range("F8").Value = 181
$ws1.Range("F9").Value = 367
$ws1.Range("F10").Value = 477
$ws1.Range("F11").Value = 519
$ws1.Range("F12").Value = 149
$ws1.Range("F13").Value = 12129
$ws1.Range("F14").Value = 5447

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 110

# --- Sheet "全部类型" (all types, aggregated) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 110
$ws4.Range("F7").Value = 7
$ws4.Range("F8").Value = 159
$ws4.Range("F10").Value = 181
$ws4.Range("F11").Value = 367
$ws4.Range("F12").Value = 477
$ws4.Range("F13").Value = 519
$ws4.Range("F14").Value = 149
$ws4.Range("F15").Value = 12129
$ws4.Range("F17").Value = 5447
